{"js": "// Update the date heading and the 25 multiplication problems in the table.\n// Cell values are set by (row, col) index to avoid any ambiguity from\n// duplicate problem text (e.g. \"35\u00d726=\" appears twice in the source).\n\nconst body = context.document.body;\n\n// 1) Update the date paragraph (first paragraph in the document body).\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[0].getRange().insertText(\"2025-07-16 Wednesday\", \"Replace\");\n\n// 2) Update the multiplication problems inside the single table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Map of (rowIndex, colIndex) -> new text, in document order.\nconst updates = [\n  // row 0\n  [0, 0, \"84\u00d796=\"],\n  [0, 1, \"43\u00d757=\"],\n  [0, 2, \"97\u00d721=\"],\n  [0, 3, \"13\u00d728=\"],\n  [0, 4, \"60\u00d792=\"],\n  // row 4\n  [4, 0, \"86\u00d749=\"],\n  [4, 1, \"77\u00d790=\"],\n  [4, 2, \"61\u00d747=\"],\n  [4, 3, \"84\u00d781=\"],\n  [4, 4, \"31\u00d753=\"],\n  // row 9\n  [9, 0, \"30\u00d788=\"],\n  [9, 1, \"13\u00d759=\"],\n  [9, 2, \"92\u00d748=\"],\n  [9, 3, \"53\u00d789=\"],\n  [9, 4, \"32\u00d713=\"],\n  // row 14\n  [14, 0, \"35\u00d726=\"],\n  [14, 1, \"82\u00d789=\"],\n  [14, 2, \"21\u00d729=\"],\n  [14, 3, \"43\u00d753=\"],\n  [14, 4, \"35\u00d780=\"],\n  // row 19\n  [19, 0, \"65\u00d781=\"],\n  [19, 1, \"25\u00d773=\"],\n  [19, 2, \"73\u00d765=\"],\n  [19, 3, \"28\u00d727=\"],\n  [19, 4, \"64\u00d713=\"],\n];\n\nfor (const [r, c, text] of updates) {\n  table.getCell(r, c).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 multiplication problems in the table.\n# Cells are addressed by (row, col) index (1-based, per Word COM convention)\n# to avoid any ambiguity from duplicate problem text (e.g. \"35\u00d726=\" appears\n# twice in the source document).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-16 Wednesday\"\n\n# 2) Update the multiplication problems inside the single table.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"84\u00d796=\"\n$t.Cell(1, 2).Range.Text = \"43\u00d757=\"\n$t.Cell(1, 3).Range.Text = \"97\u00d721=\"\n$t.Cell(1, 4).Range.Text = \"13\u00d728=\"\n$t.Cell(1, 5).Range.Text = \"60\u00d792=\"\n\n$t.Cell(5, 1).Range.Text = \"86\u00d749=\"\n$t.Cell(5, 2).Range.Text = \"77\u00d790=\"\n$t.Cell(5, 3).Range.Text = \"61\u00d747=\"\n$t.Cell(5, 4).Range.Text = \"84\u00d781=\"\n$t.Cell(5, 5).Range.Text = \"31\u00d753=\"\n\n$t.Cell(10, 1).Range.Text = \"30\u00d788=\"\n$t.Cell(10, 2).Range.Text = \"13\u00d759=\"\n$t.Cell(10, 3).Range.Text = \"92\u00d748=\"\n$t.Cell(10, 4).Range.Text = \"53\u00d789=\"\n$t.Cell(10, 5).Range.Text = \"32\u00d713=\"\n\n$t.Cell(15, 1).Range.Text = \"35\u00d726=\"\n$t.Cell(15, 2).Range.Text = \"82\u00d789=\"\n$t.Cell(15, 3).Range.Text = \"21\u00d729=\"\n$t.Cell(15, 4).Range.Text = \"43\u00d753=\"\n$t.Cell(15, 5).Range.Text = \"35\u00d780=\"\n\n$t.Cell(20, 1).Range.Text = \"65\u00d781=\"\n$t.Cell(20, 2).Range.Text = \"25\u00d773=\"\n$t.Cell(20, 3).Range.Text = \"73\u00d765=\"\n$t.Cell(20, 4).Range.Text = \"28\u00d727=\"\n$t.Cell(20, 5).Range.Text = \"64\u00d713=\"\n"}
